$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.146.03"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "2.932.58"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "591.62"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("D6").Value = "145.24"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("D9").Value = "6.97"
$ws.Range("E9").Value = "  +4.11%  "
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").Value = "33.76"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "3.417.50"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").Value = "61.238.09"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").Value = "2.928.78"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").Value = "437.69"
$ws.Range("E19").Value = "  +2.16%  "
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("D21").Value = "0.679"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("D23").Value = "81.68"
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("D26").Value = "11.85"
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "2.28"
$ws.Range("E28").Value = "  +3.75%  "
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").Value = "7.03"
$ws.Range("E30").Value = "  -2.85%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "0.111"
$ws.Range("E31").Value = "  +3.76%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "26.66"
$ws.Range("E32").Value = "  +1.03%  "
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "0.0₃0871"
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("D36").Value = "5.65"
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("D37").Value = "3.01"
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("D38").Value = "2.00"
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").Value = "0.123"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("D41").Value = "42.20"
$ws.Range("E41").Value = "  +2.21%  "
$ws.Range("D42").Value = "0.289"
$ws.Range("E42").Value = "  -2.19%  "
$ws.Range("D43").Value = "378.03"
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("D44").Value = "0.0348"
$ws.Range("E44").Value = "  -0.78%  "
$ws.Range("D45").Value = "2.693.11"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "133.28"
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("D48").Value = "23.99"
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("D51").Value = "0.125"
$ws.Range("E51").Value = "  +0.35%  "
